$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B1 value
$ws.Range("B1").Value = 25

# Delete row 2 (A2, B2) entirely
$ws.Range("A2:B2").ClearContents()

# Update selection to B1
[void]$ws.Range("B1").Select()
